$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.803.69'
$ws.Range('E2').Value = '  +2.29%  '
$ws.Range('D3').Value = '2.945.81'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'593.30"
$ws.Range('E5').Value = '  -0.30%  '
$ws.Range('D6').Value = "'147.66"
$ws.Range('E6').Value = '  +1.86%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  +1.02%  '
$ws.Range('D9').Value = '2.944.29'
$ws.Range('E9').Value = '  +0.47%  '
$ws.Range('D10').Value = "'7.28"
$ws.Range('E10').Value = '  +4.44%  '
$ws.Range('E11').Value = '  +6.92%  '
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('D13').Value = "'0.0000239"
$ws.Range('E13').Value = '  +6.18%  '
$ws.Range('D14').Value = "'32.85"
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('E15').Value = '  -0.80%  '
$ws.Range('D16').Value = '3.433.62'
$ws.Range('E16').Value = '  +0.52%  '
$ws.Range('D17').Value = '62.704.94'
$ws.Range('E17').Value = '  +2.09%  '
$ws.Range('D18').Value = "'6.70"
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('D19').Value = '2.993.61'
$ws.Range('E19').Value = '  +2.16%  '
$ws.Range('D20').Value = "'442.08"
$ws.Range('E20').Value = '  +2.24%  '
$ws.Range('D21').Value = "'13.46"
$ws.Range('E21').Value = '  -0.15%  '
$ws.Range('D22').Value = "'0.667"
$ws.Range('E22').Value = '  -1.74%  '
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('D24').Value = "'81.36"
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').Value = "'11.14"
$ws.Range('E25').Value = '  +2.11%  '
$ws.Range('D26').Value = "'2.14"
$ws.Range('E26').Value = '  -2.41%  '
$ws.Range('D27').Value = "'11.74"
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('D30').Value = "'7.24"
$ws.Range('E30').Value = '  +4.59%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = "'2.61"
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').Value = "'0.0000104"
$ws.Range('E32').Value = '  +17.83%  '
$ws.Range('D33').Value = "'26.49"
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('E34').Value = '  -0.84%  '
$ws.Range('D35').Value = "'0.999"
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('D36').Value = "'0.991"
$ws.Range('E36').Value = '  -2.15%  '
$ws.Range('D37').Value = "'3.18"
$ws.Range('E37').Value = '  +5.70%  '
$ws.Range('D38').Value = "'5.60"
$ws.Range('E38').Value = '  -0.77%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = "'49.65"
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = "'2.04"
$ws.Range('E40').Value = '  +1.70%  '
$ws.Range('D41').Value = "'8.50"
$ws.Range('E41').Value = '  -1.30%  '
$ws.Range('D42').Value = "'0.117"
$ws.Range('E42').Value = '  -5.32%  '
$ws.Range('D43').Value = "'0.281"
$ws.Range('D44').Value = "'39.76"
$ws.Range('E44').Value = '  -6.57%  '
$ws.Range('D45').Value = '2.698.72'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').Value = "'135.29"
$ws.Range('E46').Value = '  +0.66%  '
$ws.Range('D47').Value = "'0.0338"
$ws.Range('E47').Value = '  -2.49%  '
$ws.Range('D48').Value = "'363.81"
$ws.Range('E48').Value = '  -0.24%  '
$ws.Range('E50').Value = '  -0.48%  '
$ws.Range('D51').Value = "'22.89"
$ws.Range('E51').Value = '  -3.49%  '

# Reset style on quote-prefixed numeric-looking text cells so no stray
# number-format / quote-prefix styling is left behind (matches source formatting).
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
